$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 8 through 28, leaving only rows 1-7
$ws.Range("A8:A28").EntireRow.Delete()

# Replace rows 2-7 with the consolidated card text
$ws.Range("A2").Value = "('Fireball', ['{X}{R}', 'Sorcery', 'This spell costs {1} more to cast for each target beyond the first.', 'Fireball deals X damage divided evenly, rounded down, among any number of targets.'])"
$ws.Range("A3").Value = "('Mana Leak', ['{1}{U}', 'Instant', 'Counter target spell unless its controller pays {3}.'])"
$ws.Range("A4").Value = "('Oxidize', ['{G}', 'Instant', 'Destroy target artifact. It can" + [char]0x2019 + "t be regenerated.'])"
$ws.Range("A5").Value = "('Psychatog', ['{1}{U}{B}', 'Creature " + [char]0x2014 + " Atog', 'Discard a card: Psychatog gets +1/+1 until end of turn.', 'Exile two cards from your graveyard: Psychatog gets +1/+1 until end of turn.', '1/2'])"
$ws.Range("A6").Value = "('Reciprocate', ['{W}', 'Instant', 'Exile target creature that dealt damage to you this turn.'])"
$ws.Range("A7").Value = "('Terror', ['{1}{B}', 'Instant', 'Destroy target nonartifact, nonblack creature. It can" + [char]0x2019 + "t be regenerated.'])"
